$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update Notifications0008 test case (row 9) with the two additional Jira IDs / scenarios
# (set the description text first so the shared-string table keeps the same
# ordering as the saved workbook: description before the Jira id list)
$ws.Range("C9").Value = "Verify that user receives a notification when someone he is following  publishes a post||Verify that user is receiving notification when someone liked his post(aggregated notification)||Verify that user able to recevies a notification when other user commented on his post||Verify that all users receive notification when other user published a comment on post and validate notification.||Verify that all users receive notification when other user published a post and validate notification."
$ws.Range("B9").Value = "OPQA-877||OPQA-1013||OPQA-215||OPQA-1397||OPQA-1395"

# Row grew taller to accommodate the extra text
$ws.Rows.Item(9).RowHeight = 60

# Update the active cell/selection to match the saved view state
$ws.Activate()
$ws.Range("C8").Select()
